$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    # Force the literal text even when it looks like a number
    # (e.g. "315.16"), then strip the text NumberFormat back off
    # so the cell keeps its original (default) style.
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).ClearFormats()
}

$ws.Range('D2').Value = '27.311.00'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '1.833.06'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.92%  '
Set-TextValue $ws 'D5' '315.16'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('E7').Value = '  +1.79%  '
Set-TextValue $ws 'D8' '0.3688'
$ws.Range('E8').Value = '  +0.85%  '
Set-TextValue $ws 'D9' '0.07444'
$ws.Range('E9').Value = '  +1.02%  '
Set-TextValue $ws 'D10' '0.8856'
$ws.Range('E10').Value = '  +1.73%  '
Set-TextValue $ws 'D11' '20.50'
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').Value = '1.874.32'
$ws.Range('E12').Value = '  +2.48%  '
Set-TextValue $ws 'D13' '0.07329'
$ws.Range('E13').Value = '  +3.00%  '
Set-TextValue $ws 'D14' '5.436'
$ws.Range('E14').Value = '  +0.98%  '
Set-TextValue $ws 'D15' '94.01'
$ws.Range('E15').Value = '  +2.91%  '
Set-TextValue $ws 'D16' '6.565'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D20').Value = '27.533.52'
$ws.Range('E20').Value = '  +2.10%  '
Set-TextValue $ws 'D21' '14.78'
$ws.Range('E21').Value = '  +0.91%  '
Set-TextValue $ws 'D22' '5.286'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').Value = '2.095.70'
$ws.Range('E24').Value = '  +1.89%  '
Set-TextValue $ws 'D25' '1.896'
$ws.Range('E25').Value = '  +0.16%  '
Set-TextValue $ws 'D26' '152.07'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  +1.45%  '
Set-TextValue $ws 'D28' '2.144'
$ws.Range('E28').Value = '  +0.30%  '
Set-TextValue $ws 'D29' '5.230'
$ws.Range('E29').Value = '  -0.54%  '
Set-TextValue $ws 'D30' '117.11'
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('E32').Value = '  -1.24%  '
Set-TextValue $ws 'D33' '1.174'
$ws.Range('E33').Value = '  +0.76%  '
Set-TextValue $ws 'D34' '4.545'
$ws.Range('E34').Value = '  +1.45%  '
Set-TextValue $ws 'D35' '2.945'
$ws.Range('E35').Value = '  +1.66%  '
$ws.Range('E36').Value = '  +0.94%  '
Set-TextValue $ws 'D37' '1.094'
$ws.Range('E37').Value = '  -0.04%  '
Set-TextValue $ws 'D38' '0.05346'
$ws.Range('E38').Value = '  +1.12%  '
Set-TextValue $ws 'D39' '0.01953'
Set-TextValue $ws 'D40' '2.972'
$ws.Range('E40').Value = '  -0.28%  '
Set-TextValue $ws 'D41' '2.398'
$ws.Range('E41').Value = '  +3.09%  '
Set-TextValue $ws 'D42' '7.230'
$ws.Range('E42').Value = '  +1.02%  '
Set-TextValue $ws 'D43' '0.5299'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('E44').Value = '  +0.30%  '
Set-TextValue $ws 'D45' '8.476'
$ws.Range('E45').Value = '  +0.56%  '
Set-TextValue $ws 'D46' '0.4936'
$ws.Range('E46').Value = '  +1.72%  '
Set-TextValue $ws 'D47' '10.52'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('E48').Value = '  +0.90%  '
Set-TextValue $ws 'D49' '105.05'
$ws.Range('E49').Value = '  +1.69%  '
Set-TextValue $ws 'D50' '1.671'
$ws.Range('E50').Value = '  +0.64%  '
Set-TextValue $ws 'D51' '0.06299'
$ws.Range('E51').Value = '  +0.10%  '
